$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 79804
$ws.Range("B2").Value = "Dra. Maysa Teixeira"
$ws.Range("C2").Value = "Operações"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45097
$ws.Range("G2").Value = 10299.69

# Row 3
$ws.Range("A3").Value = 87854
$ws.Range("B3").Value = "Natália Vieira"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45087
$ws.Range("G3").Value = 6984.09

# Row 4
$ws.Range("A4").Value = 61727
$ws.Range("B4").Value = "Bruno Caldeira"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45084
$ws.Range("G4").Value = 5643.81

# Row 5
$ws.Range("A5").Value = 17789
$ws.Range("B5").Value = "Esther Cardoso"
$ws.Range("C5").Value = "Recursos Humanos"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45094
$ws.Range("G5").Value = 10306.19

# Row 6
$ws.Range("A6").Value = 42565
$ws.Range("B6").Value = "Luna Viana"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45084
$ws.Range("G6").Value = 5307.96

# Row 7
$ws.Range("A7").Value = 16531
$ws.Range("B7").Value = "Emanuel Barros"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45097
$ws.Range("G7").Value = 8823.129999999999

# Row 8
$ws.Range("A8").Value = 95987
$ws.Range("B8").Value = "João Vitor Nogueira"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45097
$ws.Range("G8").Value = 10818.41

# Row 9
$ws.Range("A9").Value = 70510
$ws.Range("B9").Value = "Rodrigo Novaes"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45089
$ws.Range("G9").Value = 11022.34

# Row 10
$ws.Range("A10").Value = 68228
$ws.Range("B10").Value = "Júlia Aragão"
$ws.Range("C10").Value = "Marketing"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45100
$ws.Range("G10").Value = 9280.030000000001

# Row 11
$ws.Range("A11").Value = 59885
$ws.Range("B11").Value = "Lorena Pires"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45086
$ws.Range("G11").Value = 7328.64
